# Applies the "cryptos" price/volume refresh for Thu Jun 13 10:38:02 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) sometimes holds numeric-looking text such as
# "605.30", "1.00" or "177.00". Assigning such strings directly to
# Range.Value lets Excel auto-convert them into real numbers (losing
# trailing zeros / exact text). To keep them as text - matching the
# original inline-string cells - we temporarily force the cell to Text
# number format before assigning the value, then clear the formatting
# again so the cell keeps its original (default) style.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '67.405.06'
$ws.Range("E2").Value = '  -0.67%  '

Set-TextValue $ws.Range("D3") '3.489.42'
$ws.Range("E3").Value = '  -1.51%  '

$ws.Range("E4").Value = '  -0.10%  '

Set-TextValue $ws.Range("D5") '605.30'
$ws.Range("E5").Value = '  -1.71%  '

Set-TextValue $ws.Range("D6") '150.92'
$ws.Range("E6").Value = '  -1.29%  '

Set-TextValue $ws.Range("D7") '3.486.25'
$ws.Range("E7").Value = '  -1.57%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  +1.07%  '

$ws.Range("E10").Value = '  +2.48%  '

$ws.Range("E11").Value = '  +6.33%  '

$ws.Range("E12").Value = '  +1.02%  '

$ws.Range("E13").Value = '  -1.84%  '

Set-TextValue $ws.Range("D14") '32.05'
$ws.Range("E14").Value = '  -0.39%  '

Set-TextValue $ws.Range("D15") '4.071.55'
$ws.Range("E15").Value = '  -1.72%  '

Set-TextValue $ws.Range("D16") '3.485.07'
$ws.Range("E16").Value = '  -2.00%  '

Set-TextValue $ws.Range("D17") '67.362.21'
$ws.Range("E17").Value = '  -0.45%  '

$ws.Range("E18").Value = '  -0.38%  '

Set-TextValue $ws.Range("D19") '6.49'
$ws.Range("E19").Value = '  +1.18%  '

Set-TextValue $ws.Range("D20") '15.44'
$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("E21").Value = '  +2.29%  '

Set-TextValue $ws.Range("D22") '446.28'
$ws.Range("E22").Value = '  -0.35%  '

$ws.Range("E23").Value = '  +0.13%  '

Set-TextValue $ws.Range("D24") '79.25'
$ws.Range("E24").Value = '  +2.17%  '

$ws.Range("E25").Value = '  -0.04%  '

Set-TextValue $ws.Range("D26") '3.623.45'
$ws.Range("E26").Value = '  -1.69%  '

$ws.Range("E27").Value = '  -5.10%  '

Set-TextValue $ws.Range("D28") '8.66'
$ws.Range("E28").Value = '  +0.23%  '

Set-TextValue $ws.Range("D29") '9.93'
$ws.Range("E29").Value = '  -3.47%  '

$ws.Range("E30").Value = '  -1.44%  '

$ws.Range("E31").Value = '  +1.97%  '

Set-TextValue $ws.Range("D32") '0.171'
$ws.Range("E32").Value = '  +1.21%  '

Set-TextValue $ws.Range("D33") '1.00'
$ws.Range("E33").Value = '  -0.07%  '

Set-TextValue $ws.Range("D34") '25.59'
$ws.Range("E34").Value = '  -1.50%  '

$ws.Range("E35").Value = '  -1.56%  '

Set-TextValue $ws.Range("D37") '3.477.75'
$ws.Range("E37").Value = '  -1.44%  '

Set-TextValue $ws.Range("D38") '7.99'
$ws.Range("E38").Value = '  -0.75%  '

$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("E40").Value = '  +5.14%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D41") '0.999'
$ws.Range("E41").Value = '  -0.11%  '

$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D42") '177.00'
$ws.Range("E42").Value = '  +0.09%  '

Set-TextValue $ws.Range("D43") '0.0897'
$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("E44").Value = '  -0.20%  '

Set-TextValue $ws.Range("D45") '0.891'
$ws.Range("E45").Value = '  +0.40%  '

Set-TextValue $ws.Range("D46") '30.19'
$ws.Range("E46").Value = '  +5.65%  '

$ws.Range("E47").Value = '  +2.44%  '

Set-TextValue $ws.Range("D48") '1.28'
$ws.Range("E48").Value = '  -0.19%  '

$ws.Range("E49").Value = '  -5.78%  '

Set-TextValue $ws.Range("D50") '7.60'
$ws.Range("E50").Value = '  -0.40%  '

Set-TextValue $ws.Range("D51") '0.252'
$ws.Range("E51").Value = '  -0.40%  '
